# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy the formatting used by the other header cells (G1)
# and set its text to "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cells H2:H6 = 0 (plain numeric, same as the rest of the data rows)
$ws.Range("H2:H6").Value = 0
